$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "쇼핑,백화점,마트/편의점,온라인쇼핑,소셜커머스,해외직구,홈쇼핑,SPA브랜드,아울렛,대형마트,SSM,전통시장,면세점,모든가맹점,국내외가맹점,화장품,드럭스토어"
$ws.Range("B3").Value = "문화,여가,문화/여가,영화,공연/전시,문화센터,도서,음원사이트,영화/문화,디지털구독,골프"
$ws.Range("B4").Value = "교통,고속버스,항공권,기차,대중교통,렌터카,택시,교통"
$ws.Range("B5").Value = "자동차,주유,자동차,정비,주유소,충전소,하이패스,자동차/하이패스"
$ws.Range("B6").Value = "음식,푸드,카페,패밀리레스토랑,배달앱,베이커리,점심,저녁,일반음식점,패스트푸드,카페/디저트,푸드"
$ws.Range("B7").Value = "교육,학습지,유치원,어린이집,학원,교육/육아"
$ws.Range("B8").Value = "주거,통신,주거/통신,공과금,공과금/렌탈,통신,KT,SKT,LGU+"
$ws.Range("B9").Value = "여행,항공,여행/항공,진에어,대한항공,아시아나항공,제주항공,항공마일리지,온라인 여행사,여행/숙박,공항라운지,리조트,공항,공항라운지/PP,여행사,호텔,PP,라운지키,해외이용,저가항공,해외"
$ws.Range("B10").Value = "의료,약국,병원,병원/약국"
$ws.Range("B11").Value = "금융,금융/보험,은행사,보험,보험사"
